$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels in D1/E1
$ws.Range("D1").Value = "upper"
$ws.Range("E1").Value = "lower"

# New data rows 2-11 (A,B,C,D,E,F,G,H,I,J)
$data = @(
    @(4,365,"Other and unspecified malignant neoplasm of skin",1.33,1.08,1.2,0.01,0.01,1457,"1.2[1.08, 1.33]"),
    @(4,365,"Malignant neoplasm of trachea, bronchus and lung",2.06,1.24,1.62,0,0,244,"1.62[1.24, 2.06]"),
    @(4,365,"Other malignant neoplasms of lymphoid and histiocytic tissue",2.47,1.41,1.9,0,0,205,"1.9[1.41, 2.47]"),
    @(4,365,"Malignant melanoma of skin",1.41,0.69,1,0,0,101,"1.0[0.69, 1.41]"),
    @(4,365,"Multiple myeloma and immunoproliferative neoplasms",3.56,1.26,2.15,0,0,73,"2.15[1.26, 3.56]"),
    @(4,365,"Lymphoid leukemia",3.32,0.97,1.83,0,0,47,"1.83[0.97, 3.32]"),
    @(4,365,"Myeloid leukemia",4.08,0.95,2.09,0,0,35,"2.09[0.95, 4.08]"),
    @(4,365,"Hodgkin's disease",3.93,0.85,1.99,0,0,18,"1.99[0.85, 3.93]"),
    @(4,365,"Leukemia of unspecified cell type",5.56,0.53,1.87,0,0,16,"1.87[0.53, 5.56]"),
    @(4,365,"Malignant neoplasm of larynx",5.33,0.65,2.21,0,0,15,"2.21[0.65, 5.33]")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $r++
}
